$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = "ReportStatus"
$ws.Range("I4").Font.Bold = $true

$ws.Range("I5").Value = "Authoring (Default)"

$ws.Range("I6").Select()
